$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoFCRfP")
$originallyActive = $wb.ActiveSheet

# Rename the existing "hydrogen" entry (row 24) to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add a new row 25 for "hydrogen combined cycle"
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 1

# Style A24 with a black font color and vertical-centered alignment
$a24 = $ws.Range("A24")
$a24.Font.Color = 0
$a24.VerticalAlignment = -4108

# Propagate the exact same formatting to A25 via a format-only paste so we
# don't create an extra unused style entry by re-deriving it property by
# property.
$a24.Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to match the post-edit state, then restore whichever
# sheet was active beforehand so the workbook-level active tab is unchanged.
$ws.Range("D26").Select()
$originallyActive.Activate()
